$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.559.83"
$ws.Range("E2").Value = "  +0.07%  "
$ws.Range("D3").Value = "2.645.82"
$ws.Range("E3").Value = "  +0.21%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Value = "'603.14"
$ws.Range("E5").Value = "  +2.18%  "
$ws.Range("D6").Value = "'146.56"
$ws.Range("E6").Value = "  +1.97%  "
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("E8").Value = "  +0.59%  "
$ws.Range("E9").Value = "  +1.02%  "
$ws.Range("D10").Value = "'5.59"
$ws.Range("E10").Value = "  -0.41%  "
$ws.Range("E11").Value = "  +4.57%  "
$ws.Range("E12").Value = "  -0.07%  "
$ws.Range("D13").Value = "'27.53"
$ws.Range("E13").Value = "  +0.72%  "
$ws.Range("D14").Value = "3.121.98"
$ws.Range("E14").Value = "  -0.04%  "
$ws.Range("D15").Value = "63.423.82"
$ws.Range("E15").Value = "  -0.03%  "
$ws.Range("E16").Value = "  +0.52%  "
$ws.Range("D17").Value = "2.687.72"
$ws.Range("E17").Value = "  +1.19%  "
$ws.Range("D18").Value = "'11.47"
$ws.Range("E18").Value = "  +0.55%  "
$ws.Range("D19").Value = "'4.56"
$ws.Range("E19").Value = "  +4.73%  "
$ws.Range("D20").Value = "'342.39"
$ws.Range("E20").Value = "  +0.61%  "
$ws.Range("D21").Value = "'6.93"
$ws.Range("E21").Value = "  +3.14%  "
$ws.Range("E22").Value = "  -0.10%  "
$ws.Range("D23").Value = "'5.57"
$ws.Range("E23").Value = "  -3.24%  "
$ws.Range("D24").Value = "'66.69"
$ws.Range("E24").Value = "  -0.97%  "
$ws.Range("D25").Value = "'1.69"
$ws.Range("E25").Value = "  +0.07%  "
$ws.Range("D26").Value = "'9.07"
$ws.Range("E26").Value = "  +7.23%  "
$ws.Range("D27").Value = "'568.91"
$ws.Range("E27").Value = "  +3.15%  "
$ws.Range("D28").Value = "'1.55"
$ws.Range("E28").Value = "  +0.91%  "
$ws.Range("E29").Value = "  -1.41%  "
$ws.Range("E32").Value = "  +3.07%  "
$ws.Range("E33").Value = "  -3.45%  "
$ws.Range("D34").Value = "0.0₃0816"
$ws.Range("E34").Value = "  +1.28%  "
$ws.Range("D35").Value = "'5.18"
$ws.Range("E35").Value = "  +5.11%  "
$ws.Range("D36").Value = "'168.44"
$ws.Range("E36").Value = "  -3.71%  "
$ws.Range("D37").Value = "'0.408"
$ws.Range("E37").Value = "  +1.30%  "
$ws.Range("D38").Value = "'0.999"
$ws.Range("E38").Value = "  -0.11%  "
$ws.Range("D39").Value = "'1.93"
$ws.Range("E39").Value = "  +6.10%  "
$ws.Range("E40").Value = "  +0.15%  "
$ws.Range("E41").Value = "  -0.10%  "
$ws.Range("D42").Value = "'168.64"
$ws.Range("E42").Value = "  -0.97%  "
$ws.Range("D43").Value = "'3.78"
$ws.Range("E43").Value = "  +1.23%  "
$ws.Range("D44").Value = "'22.20"
$ws.Range("E44").Value = "  -0.46%  "
$ws.Range("E45").Value = "  +2.94%  "
$ws.Range("D46").Value = "'0.631"
$ws.Range("E46").Value = "  +0.13%  "
$ws.Range("D47").Value = "'0.0246"
$ws.Range("E47").Value = "  +3.45%  "
$ws.Range("D48").Value = "'0.0961"
$ws.Range("E48").Value = "  +0.11%  "
$ws.Range("D49").Value = "'18.82"
$ws.Range("E49").Value = "  +0.55%  "
$ws.Range("D50").Value = "'1.88"
$ws.Range("E50").Value = "  +9.73%  "
$ws.Range("B30").Value = "Aptos"
$ws.Range("C30").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D30").Value = "'8.00"
$ws.Range("E30").Value = "  +2.70%  "
$ws.Range("B31").Value = "Binance-PegBSC-USD"
$ws.Range("C31").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D31").Value = "'0.999"
$ws.Range("E31").Value = "  -0.16%  "
$ws.Range("B51").Value = "WhiteBITCoin"
$ws.Range("C51").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D51").Value = "'11.27"
$ws.Range("E51").Value = "  -0.77%  "
